$wb = $excel.ActiveWorkbook

# Add a "Comments" column (column E) header to the history sheets that
# currently have data in columns A:D, and leave cell E1 selected on each
# of those sheets.
$sheetNames = @("Withdraw History", "Deposit History", "Transfer History", "Absolute History")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Select()
    $ws.Range("E1").Value = "Comments"
    $ws.Range("E1").Select()
}

# Make "Withdraw History" the active/selected sheet (instead of "Absolute History").
$wsActive = $wb.Worksheets.Item("Withdraw History")
$wsActive.Select()
$wsActive.Range("E1").Select()
